$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (F column = 想去人数 / "people interested" counts)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 16
$wsExhibit.Range("F9").Value = 6860
$wsExhibit.Range("F14").Value = 127
$wsExhibit.Range("F16").Value = 4369
$wsExhibit.Range("F20").Value = 4356
$wsExhibit.Range("F31").Value = 72
$wsExhibit.Range("F32").Value = 7912
$wsExhibit.Range("F35").Value = 661
$wsExhibit.Range("F39").Value = 1592
$wsExhibit.Range("F41").Value = 915
$wsExhibit.Range("F43").Value = 3981

# Sheet "全部类型" updates (same events, same F column data)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 16
$wsAll.Range("F12").Value = 6860
$wsAll.Range("F17").Value = 127
$wsAll.Range("F19").Value = 4369
$wsAll.Range("F22").Value = 4356
$wsAll.Range("F31").Value = 72
$wsAll.Range("F33").Value = 7912
$wsAll.Range("F36").Value = 661
$wsAll.Range("F39").Value = 1592
$wsAll.Range("F41").Value = 915
$wsAll.Range("F43").Value = 3981
